$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper note: cells whose new value looks like a plain decimal number
# (single '.' and only digits) need NumberFormat forced to Text ("@")
# before assignment, otherwise Excel's COM layer auto-converts the
# string into a numeric value and trailing zeros / formatting are lost.
# Cells whose values are not number-like (multiple dots, letters, URLs,
# percent strings padded with spaces, etc.) are safe to set directly.

# Row 2 - Bitcoin
$ws.Range("D2").Value = "28.357.01"
$ws.Range("E2").Value = "  -0.48%  "

# Row 3 - Ethereum
$ws.Range("D3").Value = "1.570.18"
$ws.Range("E3").Value = "  +0.53%  "

# Row 4 - TetherUSD
$ws.Range("E4").Value = "  +0.09%  "

# Row 5 - BNB
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "211.34"
$ws.Range("E5").Value = "  -0.11%  "

# Row 6 - XRP
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.492"
$ws.Range("E6").Value = "  -0.17%  "

# Row 7 - USDC
$ws.Range("E7").Value = "  +0.08%  "

# Row 8 - OKB
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "44.26"
$ws.Range("E8").Value = "  -3.73%  "

# Row 9 - Solana
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "23.81"
$ws.Range("E9").Value = "  -1.04%  "

# Row 10 - Cardano
$ws.Range("E10").Value = "  -0.64%  "

# Row 12 - TRON
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.0896"
$ws.Range("E12").Value = "  +1.59%  "

# Row 13 - WrappedliquidstakedEther2.0
$ws.Range("D13").Value = "1.794.39"
$ws.Range("E13").Value = "  +0.58%  "

# Row 14 - WrappedEther
$ws.Range("D14").Value = "1.587.79"
$ws.Range("E14").Value = "  +1.49%  "

# Row 15 - Polkadot
$ws.Range("E15").Value = "  -0.12%  "

# Row 16 - WrappedBTC
$ws.Range("D16").Value = "28.353.15"
$ws.Range("E16").Value = "  -0.48%  "

# Row 17 - Polygon
$ws.Range("E17").Value = "  -1.02%  "

# Row 18 - Litecoin
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "61.29"
$ws.Range("E18").Value = "  -1.02%  "

# Row 19 - BitcoinCash
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "227.91"
$ws.Range("E19").Value = "  +0.51%  "

# Row 20 - Chainlink
$ws.Range("E20").Value = "  +0.97%  "

# Row 21 - ShibaInu
$ws.Range("E21").Value = "  -1.35%  "

# Row 22 - Dai
$ws.Range("E22").Value = "  +0.07%  "

# Row 23 - Uniswap
$ws.Range("E23").Value = "  +1.86%  "

# Row 24 - Avalanche
$ws.Range("E24").Value = "  -1.75%  "

# Row 25 - Toncoin
$ws.Range("E25").Value = "  -0.77%  "

# Row 26 - Monero
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "150.96"
$ws.Range("E26").Value = "  +0.75%  "

# Row 27 - EthereumClassic
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "14.91"
$ws.Range("E27").Value = "  -0.31%  "

# Row 28 - Stellar
$ws.Range("E28").Value = "  -0.39%  "

# Row 29 - Cosmos
$ws.Range("E29").Value = "  -1.29%  "

# Row 30 - BinanceUSD
$ws.Range("E30").Value = "  +0.08%  "

# Row 31 - Hedera
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.0481"
$ws.Range("E31").Value = "  +3.66%  "

# Row 32 - PancakeSwap
$ws.Range("E32").Value = "  -2.56%  "

# Row 33 - Filecoin
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.18"
$ws.Range("E33").Value = "  -0.57%  "

# Row 34 - InternetComputer(DFINITY)
$ws.Range("E34").Value = "  -1.63%  "

# Row 35 - Maker
$ws.Range("D35").Value = "1.380.20"
$ws.Range("E35").Value = "  -0.98%  "

# Row 36 - TrustWalletToken
$ws.Range("E36").Value = "  +2.29%  "

# Row 37 - LidoDAOToken
$ws.Range("E37").Value = "  -2.64%  "

# Row 38 - HuobiToken
$ws.Range("E38").Value = "  -0.28%  "

# Row 39 - MXToken
$ws.Range("E39").Value = "  +1.75%  "

# Row 40 - VeChain
$ws.Range("E40").Value = "  -1.65%  "

# Row 41 - ImmutableX
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.520"
$ws.Range("E41").Value = "  -2.64%  "

# Row 42 - now RenderToken (was PaxDollar)
$ws.Range("B42").Value = "RenderToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.91"
$ws.Range("E42").Value = "  +3.19%  "

# Row 43 - now PaxDollar (was RenderToken)
$ws.Range("B43").Value = "PaxDollar"
$ws.Range("C43").Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.00"
$ws.Range("E43").Value = "  +0.04%  "

# Row 45 - Kaspa
$ws.Range("E45").Value = "  -0.79%  "

# Row 46 - FraxShare
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "5.33"
$ws.Range("E46").Value = "  -3.60%  "

# Row 47 - Aave
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "62.30"

# Row 48 - WEMIXToken
$ws.Range("E48").Value = "  -6.18%  "

# Row 49 - RocketPoolETH
$ws.Range("D49").Value = "1.707.09"
$ws.Range("E49").Value = "  +0.55%  "

# Row 50 - mCoin
$ws.Range("E50").Value = "  +1.95%  "

# Row 51 - Quant
$ws.Range("E51").Value = "  -0.74%  "
